$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate so the rebuilt layout matches exactly.
$ws.Cells.Clear()

# Row 2: "1min" label
$ws.Range("A2").Value = "1min"

# Row 3: column group headers
$ws.Range("C3").Value = "Normal"
$ws.Range("F3").Value = "Mutant"

# Row 4: date header row (moved down from the old row 1)
$ws.Range("B4:G4").NumberFormat = "m/d/yy"
$ws.Range("B4").Value = 42648
$ws.Range("C4").Value = 42655
$ws.Range("D4").Value = 42662
$ws.Range("E4").Value = 42648
$ws.Range("F4").Value = 42655
$ws.Range("G4").Value = 42662

# Row 5: B6 data
$ws.Range("A5").Value = "B6"
$ws.Range("B5").Value = 146.6
$ws.Range("C5").Value = 138.6
$ws.Range("D5").Value = 155.6
$ws.Range("E5").Value = 166
$ws.Range("F5").Value = 179.3
$ws.Range("G5").Value = 186.9

# Row 6: BTBR data
$ws.Range("A6").Value = "BTBR"
$ws.Range("B6").Value = 245.7
$ws.Range("C6").Value = 240
$ws.Range("D6").Value = 243.1
$ws.Range("E6").Value = 177.8
$ws.Range("F6").Value = 171.6
$ws.Range("G6").Value = 188.1

# Row 8: "5min" label
$ws.Range("A8").Value = "5min"

# Row 9: column group headers (second block)
$ws.Range("C9").Value = "Normal"
$ws.Range("F9").Value = "Mutant"

# Row 10: date header row (new, mirrors row 4)
$ws.Range("B10:G10").NumberFormat = "m/d/yy"
$ws.Range("B10").Value = 42648
$ws.Range("C10").Value = 42655
$ws.Range("D10").Value = 42662
$ws.Range("E10").Value = 42648
$ws.Range("F10").Value = 42655
$ws.Range("G10").Value = 42662

# Row 11: B6 data (second block) - keeps its custom font style
$ws.Range("A11").Value = "B6"
$ws.Range("A11").Font.Color = 0
$ws.Range("B11").Value = 333.6
$ws.Range("C11").Value = 353.6
$ws.Range("D11").Value = 408.8
$ws.Range("E11").Value = 450.6
$ws.Range("F11").Value = 474.4
$ws.Range("G11").Value = 423.8

# Row 12: BTBR data (second block) - keeps its custom font style
$ws.Range("A12").Value = "BTBR"
$ws.Range("A12").Font.Color = 0
$ws.Range("B12").Value = 514.4
$ws.Range("C12").Value = 610.6
$ws.Range("D12").Value = 597.9
$ws.Range("E12").Value = 412.1
$ws.Range("F12").Value = 447.4
$ws.Range("G12").Value = 446.5

# Restore the selected cell shown when the sheet was last saved.
$ws.Range("C8").Select()
